$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit text format
# applied first, otherwise Excel auto-converts the text into a floating point number
# (losing formatting / precision), which the source data never intends since the
# "Price" column always holds plain text in this workbook.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D14", "D15", "D21", "D22", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D34", "D36", "D45", "D49", "D50", "D51")
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range('D2').Value = '41.671.08'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.470.38'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '316.53'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').Value = '92.65'
$ws.Range('E6').Value = '  -0.56%  '
$ws.Range('E7').Value = '  +1.83%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  +2.77%  '
$ws.Range('D10').Value = '0.0877'
$ws.Range('E10').Value = '  +11.67%  '
$ws.Range('D11').Value = '32.85'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('E12').Value = '  -0.54%  '
$ws.Range('D13').Value = '2.852.38'
$ws.Range('D14').Value = '6.90'
$ws.Range('E14').Value = '  +0.43%  '
$ws.Range('D15').Value = '15.75'
$ws.Range('E15').Value = '  -2.12%  '
$ws.Range('D16').Value = '2.460.27'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('E17').Value = '  +3.36%  '
$ws.Range('D18').Value = '41.649.82'
$ws.Range('E18').Value = '  +0.25%  '
$ws.Range('D19').Value = '0.0₃0970'
$ws.Range('E19').Value = '  +4.42%  '
$ws.Range('E20').Value = '  +2.09%  '
$ws.Range('D21').Value = '71.35'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('D22').Value = '11.43'
$ws.Range('E22').Value = '  +2.03%  '
$ws.Range('D23').Value = '238.75'
$ws.Range('E23').Value = '  +0.95%  '
$ws.Range('D24').Value = '2.73'
$ws.Range('E24').Value = '  +0.32%  '
$ws.Range('D25').Value = '1.91'
$ws.Range('E25').Value = '  -1.13%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').Value = '24.69'
$ws.Range('E27').Value = '  -2.20%  '
$ws.Range('E28').Value = '  +2.35%  '
$ws.Range('D29').Value = '9.83'
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('D30').Value = '35.44'
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('D31').Value = '155.96'
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').Value = '  +0.76%  '
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('D34').Value = '0.0764'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E35').Value = '  +2.14%  '
$ws.Range('D36').Value = '17.60'
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('E37').Value = '  -2.23%  '
$ws.Range('E38').Value = '  +1.20%  '
$ws.Range('E39').Value = '  -2.25%  '
$ws.Range('E40').Value = '  -3.06%  '
$ws.Range('E41').Value = '  -2.79%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').Value = '1.967.85'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').Value = '18.90'
$ws.Range('E45').Value = '  -5.32%  '
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('E47').Value = '  +2.05%  '
$ws.Range('D48').Value = '2.706.60'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').Value = '97.31'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('D50').Value = '66.83'
$ws.Range('E50').Value = '  -1.51%  '
$ws.Range('D51').Value = '52.57'
$ws.Range('E51').Value = '  +3.82%  '
